# REVER_DailyTracker: add Oct 6-16 2020 entries to the "Ram" sheet,
# update the active-cell selection, and refresh workbook view metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ram")
$wsV = $wb.Worksheets.Item("Vijay")
$wsN = $wb.Worksheets.Item("Nirmal")

# Style-source cells (already present in the workbook) used purely to clone
# formatting via Copy/PasteSpecial(xlPasteFormats) so we don't mint new
# cellXfs entries for formats that already exist.
$srcA = $ws.Range("A4")     # plain bordered cell
$srcB = $ws.Range("B4")     # date, custom format, border
$srcC = $ws.Range("C4")     # left aligned text, border
$srcDTask = $ws.Range("D7") # left aligned wrap text, border
$srcE = $wsV.Range("E2")    # percent, centered, border
$srcG = $ws.Range("G7")     # plain cell, border 3
$srcFWip = $ws.Range("B23") # fill 3 (WIP), centered
$srcFDone = $wsN.Range("F6")# fill 5 (Completed), border 1 (no alignment yet)
$srcDOff = $ws.Range("D4")  # centered/bold-red "Holiday/Week off" style
$srcEOff = $ws.Range("E4")  # off-day style, fill 7
$srcFOff = $ws.Range("F4")  # off-day style, fill 7

function Set-RowCommon($r) {
    $srcA.Copy(); $ws.Range("A$r").PasteSpecial(-4122)
    $srcB.Copy(); $ws.Range("B$r").PasteSpecial(-4122)
    $srcC.Copy(); $ws.Range("C$r").PasteSpecial(-4122)
    $srcG.Copy(); $ws.Range("G$r").PasteSpecial(-4122)
}

function Set-TaskRow($r, $no, $date, $app, $task, $pct, $status) {
    Set-RowCommon $r
    $srcDTask.Copy(); $ws.Range("D$r").PasteSpecial(-4122)
    $srcE.Copy(); $ws.Range("E$r").PasteSpecial(-4122)
    if ($status -eq "Completed") {
        $srcFDone.Copy(); $ws.Range("F$r").PasteSpecial(-4122)
        $ws.Range("F$r").HorizontalAlignment = -4108
    } else {
        $srcFWip.Copy(); $ws.Range("F$r").PasteSpecial(-4122)
    }

    $ws.Range("A$r").Value2 = $no
    $ws.Range("B$r").Value2 = $date
    $ws.Range("C$r").Value2 = $app
    $ws.Range("D$r").Value2 = $task
    if ($null -ne $pct) {
        $ws.Range("E$r").Value2 = $pct
    }
    $ws.Range("F$r").Value2 = $status
}

function Set-OffRow($r, $no, $date, $label) {
    Set-RowCommon $r
    $srcDOff.Copy(); $ws.Range("D$r").PasteSpecial(-4122)
    $srcEOff.Copy(); $ws.Range("E$r").PasteSpecial(-4122)
    $srcFOff.Copy(); $ws.Range("F$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value2 = $no
    $ws.Range("B$r").Value2 = $date
    $ws.Range("D$r").Value2 = $label
}

Set-TaskRow 8  6  44110 "Soniya"   "creating Setup file"                            1    "Completed"
Set-TaskRow 9  7  44111 "QMVAR"    "Layout Designing"                               $null "WIP"
Set-TaskRow 10 8  44112 "Emplogin" "Adding assests Detatils"                        0.5  "WIP"
Set-TaskRow 11 9  44113 "Emplogin" "Adding assests Detatils"                        1    "Completed"
Set-OffRow  12 10 44114 "Week off"
Set-OffRow  13 11 44115 "Week off"
Set-TaskRow 14 12 44116 "Emplogin" "Adding assests Detatils and responisble view"   1    "Completed"
Set-TaskRow 15 13 44117 "QMVAR"    "Adding assests Detatils"                        1    "WIP"
Set-TaskRow 16 14 44118 "QMVAR"    "stored management added in Monthly target"      1    "Completed"
Set-TaskRow 17 15 44119 "QMVAR"    "adding layout in analysis file upload"          0.7  "WIP"
Set-TaskRow 18 16 44120 "QMVAR"    "adding layout in analysis file Export"          0.5  "WIP"

$ws.Range("F23").Select()

Write-Output "edit complete"
